$d = $word.ActiveDocument

# Locate the "9 dict { } (dicionários)" heading paragraph; the new content
# belongs in the (currently empty) paragraph right after it.
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "dict.*\(dicion") {
        $targetPara = $p
        break
    }
}

if ($targetPara -eq $null) {
    throw "could not locate target heading paragraph"
}

$insertionPoint = $targetPara.Next()
$r = $insertionPoint.Range

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Métodos úteis dos dicionários em Python</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>len</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> - quantas chaves</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>keys</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> - iterável com as chaves</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:lastRenderedPageBreak/><w:t>values</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> - iterável com os valores</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>items</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> - iterável com chaves e valores</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>setdefault</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> - adiciona valor se a chave não existe</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>copy</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> - retorna uma cópia rasa (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>shallow</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>copy</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>get</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> - obtém uma chave</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>pop</w:t></w:r><w:r><w:t xml:space="preserve"> - Apaga um item com a chave especificada (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>del</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>popitem</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> - Apaga o último item adicionado</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>update</w:t></w:r><w:r><w:t xml:space="preserve"> - Atualiza um dicionário com outro</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$r.InsertXML($xml)

Write-Output "done; paragraph count = $($d.Paragraphs.Count)"
